# Task_小游戏任务表.xlsx - "update config table 0.9.0"
#
# The header row (row 3, the "field name" row) is being renamed/retargeted:
#   configPrefabGuid            -> questObjectGuid
#   任务备注                     -> 名称
#   事件完成次数                  -> 子项目数
#   可否重复完成                  -> 可重复性
#   任务完成奖励...(multiline)    -> 完成奖励
#   小游戏配置预制体 Guid          -> Quest 物体 Guid

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value2 = "questObjectGuid"

$ws.Range("B3").Value2 = "名称"
$ws.Range("C3").Value2 = "子项目数"
$ws.Range("D3").Value2 = "可重复性"
$ws.Range("E3").Value2 = "完成奖励"
$ws.Range("F3").Value2 = "Quest 物体 Guid"

Write-Host "Updated header labels for Task_小游戏任务表"
